# Generate Report for Handback
# Row 7 ("62a80142-a66a-4fca-822b-e35e9a7e7248") in both the zh-cn and de-de
# sheets now has a completed handback: a Latest Target File hyperlink, a
# Latest Handback File name, a Latest Handback DateTime, and an Error
# Detail describing that the handback is based on a stale commit.

$wb = $excel.ActiveWorkbook

$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/46504538b033a09db3ff1d5607cc79512efc77c5/e2e/62a80142-a66a-4fca-822b-e35e9a7e7248.md"
$hyperlinkDisplay = "62a80142-a66a-4fca-822b-e35e9a7e7248.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5ba0830646490da0ead2be7d1dbcf911e8e43b7c/e2e/62a80142-a66a-4fca-822b-e35e9a7e7248.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/46504538b033a09db3ff1d5607cc79512efc77c5/e2e/62a80142-a66a-4fca-822b-e35e9a7e7248.md."

function Update-HandbackRow($ws, [string]$handbackFile, [string]$handbackDateTime) {
    # I7 - Latest Target File: becomes a hyperlink to the source markdown,
    # mirroring the style already used for A7.
    $targetCell = $ws.Range("I7")
    $targetCell.Value = $hyperlinkDisplay
    $ws.Hyperlinks.Add($targetCell, $hyperlinkUrl, [Type]::Missing, [Type]::Missing, $hyperlinkDisplay)
    $targetCell.Font.Underline = 2
    $targetCell.Font.Color = 15570276

    # J7 - Latest Handback File
    $ws.Range("J7").Value = $handbackFile

    # K7 - Latest Handback DateTime
    $ws.Range("K7").Value = $handbackDateTime

    # P7 - Error Detail
    $ws.Range("P7").Value = $errorDetail
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow $wsZhCn "62a80142-a66a-4fca-822b-e35e9a7e7248.2db139deef9d6ec13b988ad65ff23fb940fc0a67.zh-cn.xlf" "2016-08-21 01:00:53"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow $wsDeDe "62a80142-a66a-4fca-822b-e35e9a7e7248.2db139deef9d6ec13b988ad65ff23fb940fc0a67.de-de.xlf" "2016-08-21 01:00:59"
